# running_time.xlsx update: append a second "MARS Tool Output" block
# (rows 18-32) that mirrors the existing block in rows 2-16, with new
# data values, onto Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Helper: copy ONLY the formatting of a single source cell onto a
# target cell (xlPasteFormats = -4122) so we reuse the exact existing
# style index instead of minting new ones in styles.xml.
# ---------------------------------------------------------------------
function Copy-Format($srcAddr, $dstAddr) {
    $ws.Range($srcAddr).Copy() | Out-Null
    $ws.Range($dstAddr).PasteSpecial(-4122) | Out-Null
}

# Row 18 mirrors row 2 (bold 14pt banner row: title + "Calulations" label)
Copy-Format "A2" "A18"
Copy-Format "B2" "B18"
Copy-Format "C2" "C18"
Copy-Format "D2" "D18"
$ws.Range("A18").Value = "MARS Tool Output"
$ws.Range("D18").Value = "Calulations"

# Row 20 mirrors row 4 (bold sub-heading)
Copy-Format "A4" "A20"
$ws.Range("A20").Value = "Instruction Statistics Tool"

# Row 21 mirrors row 5 (bottom-bordered column headers)
Copy-Format "A5" "A21"
Copy-Format "B5" "B21"
Copy-Format "D5" "D21"
Copy-Format "E5" "E21"
Copy-Format "F5" "F21"
$ws.Range("A21").Value = "Instruction type"
$ws.Range("B21").Value = "Count"
$ws.Range("D21").Value = "Adjusted count"
$ws.Range("E21").Value = "CPI"
$ws.Range("F21").Value = "Total cycles"

# Row 22 mirrors row 6 (ALU)
$ws.Range("A22").Value = "ALU"
$ws.Range("B22").Value = 3721
$ws.Range("D22").Formula = "=B22"
$ws.Range("E22").Value = 1
$ws.Range("F22").Formula = "=D22*E22"

# Row 23 mirrors row 7 (Jump)
$ws.Range("A23").Value = "Jump"
$ws.Range("B23").Value = 298
$ws.Range("D23").Formula = "=B23"
$ws.Range("E23").Value = 1
$ws.Range("F23").Formula = "=D23*E23"

# Row 24 mirrors row 8 (Branch)
$ws.Range("A24").Value = "Branch"
$ws.Range("B24").Value = 967
$ws.Range("D24").Formula = "=B24"
$ws.Range("E24").Value = 2
$ws.Range("F24").Formula = "=D24*E24"

# Row 25 mirrors row 9 (Memory) - no D/E/F values
$ws.Range("A25").Value = "Memory"
$ws.Range("B25").Value = 621

# Row 26 mirrors row 10 (Other, with the cross-referencing D formula)
$ws.Range("A26").Value = "Other"
$ws.Range("B26").Value = 1483
$ws.Range("D26").Formula = "=B26-(B30+B31-B25)"
$ws.Range("E26").Value = 5
$ws.Range("F26").Formula = "=D26*E26"

# Row 28 mirrors row 12 (bold sub-heading)
Copy-Format "A12" "A28"
$ws.Range("A28").Value = "Data Cache Simulation Tool"

# Row 29 mirrors row 13 (bottom-bordered column headers)
Copy-Format "A13" "A29"
Copy-Format "B13" "B29"
$ws.Range("A29").Value = "Access"
$ws.Range("B29").Value = "Count"

# Row 30 mirrors row 14 (Cache hit)
$ws.Range("A30").Value = "Cache hit"
$ws.Range("B30").Value = 362
$ws.Range("D30").Formula = "=B30"
$ws.Range("E30").Value = 2
$ws.Range("F30").Formula = "=D30*E30"

# Row 31 mirrors row 15 (Cache miss)
$ws.Range("A31").Value = "Cache miss"
$ws.Range("B31").Value = 351
$ws.Range("D31").Formula = "=B31"
$ws.Range("E31").Value = 40
$ws.Range("F31").Formula = "=D31*E31"

# Row 32 mirrors row 16 (bold 14pt + top border total row)
Copy-Format "F16" "F32"
$ws.Range("F32").Formula = "=SUM(F22:F31)"

# Row heights for the two banner rows (match rows 2 & 16, which are 18.75)
$ws.Rows.Item(18).RowHeight = 18.75
$ws.Rows.Item(32).RowHeight = 18.75

# View: scrolled down to show the new block, selection sitting on B32
$ws.Application.ActiveWindow.ScrollRow = 13
$ws.Range("B32").Select() | Out-Null
